$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values per diff ---
$ws.Range("O2").Value = 1.3
$ws.Range("N3").Value = 1.3
$ws.Range("P3").Value = 1.3
$ws.Range("R3").Value = 1.21
$ws.Range("T3").Value = 1.04
$ws.Range("U3").Value = 1.04
$ws.Range("N4").Value = 1.17
$ws.Range("P4").Value = 1.17
$ws.Range("R4").Value = 1.1
$ws.Range("S4").Value = 1.05
$ws.Range("F5").Value = 6.6
$ws.Range("I5").Value = 1.56
$ws.Range("S5").Value = 2.6
$ws.Range("Q6").Value = 1.64
$ws.Range("T6").Value = 1.04
$ws.Range("U6").Value = 1.04
$ws.Range("F7").Value = 1.72
$ws.Range("G7").Value = 1.92
$ws.Range("J7").Value = 3.6
$ws.Range("N7").Value = 2.78
$ws.Range("T7").Value = 1.04
$ws.Range("U7").Value = 1.04
$ws.Range("W7").Value = 2.08
$ws.Range("S8").Value = 3.35
$ws.Range("T8").Value = 2.22
$ws.Range("AN8").Value = 7
$ws.Range("G9").Value = 2.52
$ws.Range("L9").Value = 1.01
$ws.Range("M9").Value = 1.01
$ws.Range("N9").Value = 2.78
$ws.Range("O9").Value = 1.01
$ws.Range("R9").Value = 1.25
$ws.Range("S9").Value = 3.2
$ws.Range("T9").Value = 1.04
$ws.Range("U9").Value = 1.04
$ws.Range("V9").Value = 1.36
$ws.Range("W9").Value = 1.65
$ws.Range("X9").Value = 990
$ws.Range("Y9").Value = 990
$ws.Range("Z9").Value = 1000
$ws.Range("AA9").Value = 1000
$ws.Range("AB9").Value = 990
$ws.Range("AC9").Value = 990
$ws.Range("AD9").Value = 990
$ws.Range("AE9").Value = 1000
$ws.Range("AF9").Value = 1000
$ws.Range("AG9").Value = 990
$ws.Range("AH9").Value = 990
$ws.Range("AI9").Value = 1000
$ws.Range("AJ9").Value = 1000
$ws.Range("AK9").Value = 1000
$ws.Range("AL9").Value = 1000
$ws.Range("AM9").Value = 1000
$ws.Range("AN9").Value = 1000
$ws.Range("AO9").Value = 1000
$ws.Range("F10").Value = 3.25
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 1.92
$ws.Range("I10").Value = 2.48
$ws.Range("J10").Value = 3.05
$ws.Range("K10").Value = 5.8
$ws.Range("L10").Value = 1.01
$ws.Range("M10").Value = 1.01
$ws.Range("N10").Value = 2.66
$ws.Range("O10").Value = 1.01
$ws.Range("P10").Value = 1.63
$ws.Range("Q10").Value = 1.94
$ws.Range("R10").Value = 1.22
$ws.Range("S10").Value = 3.4
$ws.Range("T10").Value = 1.04
$ws.Range("U10").Value = 1.04
$ws.Range("V10").Value = 1.67
$ws.Range("W10").Value = 1.25
$ws.Range("X10").Value = 990
$ws.Range("Y10").Value = 990
$ws.Range("Z10").Value = 1000
$ws.Range("AA10").Value = 1000
$ws.Range("AB10").Value = 990
$ws.Range("AC10").Value = 990
$ws.Range("AD10").Value = 990
$ws.Range("AE10").Value = 1000
$ws.Range("AF10").Value = 1000
$ws.Range("AG10").Value = 990
$ws.Range("AH10").Value = 990
$ws.Range("AI10").Value = 1000
$ws.Range("AJ10").Value = 1000
$ws.Range("AK10").Value = 1000
$ws.Range("AL10").Value = 1000
$ws.Range("AM10").Value = 1000
$ws.Range("AN10").Value = 1000
$ws.Range("AO10").Value = 1000
$ws.Range("H11").Value = 1.59
$ws.Range("I11").Value = 1.61
$ws.Range("J11").Value = 4.2
$ws.Range("K11").Value = 4.3
$ws.Range("L11").Value = 1.43
$ws.Range("O11").Value = 1.38
$ws.Range("V11").Value = 2.64
$ws.Range("W11").Value = 1.15
$ws.Range("X11").Value = 12.5
$ws.Range("Z11").Value = 8.199999999999999
$ws.Range("AA11").Value = 14.5
$ws.Range("AC11").Value = 9.199999999999999
$ws.Range("AD11").Value = 10
$ws.Range("AE11").Value = 19.5
$ws.Range("AF11").Value = 60
$ws.Range("AG11").Value = 28
$ws.Range("AI11").Value = 50
$ws.Range("AJ11").Value = 250
$ws.Range("AK11").Value = 130
$ws.Range("AL11").Value = 140
$ws.Range("AM11").Value = 190
$ws.Range("AN11").Value = 210
$ws.Range("F12").Value = 2.8
$ws.Range("G12").Value = 2.82
$ws.Range("H12").Value = 2.9
$ws.Range("I12").Value = 2.92
$ws.Range("L12").Value = 1.49
$ws.Range("M12").Value = 1.1
$ws.Range("P12").Value = 1.75
$ws.Range("Q12").Value = 2.3
$ws.Range("T12").Value = 1.94
$ws.Range("V12").Value = 1.52
$ws.Range("W12").Value = 1.55
$ws.Range("Z12").Value = 17.5
$ws.Range("AA12").Value = 46
$ws.Range("AB12").Value = 9.800000000000001
$ws.Range("AC12").Value = 7.2
$ws.Range("AD12").Value = 12.5
$ws.Range("AH12").Value = 19
$ws.Range("AI12").Value = 55
$ws.Range("AK12").Value = 34
$ws.Range("AO12").Value = 38
$ws.Range("K14").Value = 3.2
$ws.Range("P16").Value = 1.91
$ws.Range("Q16").Value = 2.08
$ws.Range("T16").Value = 1.84
$ws.Range("AG16").Value = 11
$ws.Range("AI16").Value = 60
$ws.Range("AJ16").Value = 32
$ws.Range("F17").Value = 1.77
$ws.Range("G17").Value = 1.79
$ws.Range("I17").Value = 5.9
$ws.Range("AC17").Value = 8.6
$ws.Range("J18").Value = 3.3
$ws.Range("H19").Value = 2.08
$ws.Range("I19").Value = 2.12
$ws.Range("P19").Value = 2.18
$ws.Range("T19").Value = 1.71
$ws.Range("U19").Value = 2.32
$ws.Range("H20").Value = 23
$ws.Range("I20").Value = 30
$ws.Range("F19").Value = 3.95
$ws.Range("G19").Value = 4

# --- Add new row 21 (Honduras Liga Nacional) ---
$ws.Range("A21").Value = "Honduras Liga Nacional"
$ws.Range("B21").Value = "'2025-12-28"
$ws.Range("C21").Value = "'18:00:00"
$ws.Range("D21").Value = "CD Marathon"
$ws.Range("E21").Value = "Olancho"
$ws.Range("F21").Value = 1.04
$ws.Range("G21").Value = 1000
$ws.Range("H21").Value = 1.04
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 1.01
$ws.Range("K21").Value = 950
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 1.24
$ws.Range("Q21").Value = 1.01
$ws.Range("R21").Value = 0
$ws.Range("S21").Value = 0
$ws.Range("T21").Value = 0
$ws.Range("U21").Value = 0
$ws.Range("V21").Value = 0
$ws.Range("W21").Value = 0
$ws.Range("X21").Value = 0
$ws.Range("Y21").Value = 0
$ws.Range("Z21").Value = 0
$ws.Range("AA21").Value = 0
$ws.Range("AB21").Value = 0
$ws.Range("AC21").Value = 0
$ws.Range("AD21").Value = 0
$ws.Range("AE21").Value = 0
$ws.Range("AF21").Value = 0
$ws.Range("AG21").Value = 0
$ws.Range("AH21").Value = 0
$ws.Range("AI21").Value = 0
$ws.Range("AJ21").Value = 0
$ws.Range("AK21").Value = 0
$ws.Range("AL21").Value = 0
$ws.Range("AM21").Value = 0
$ws.Range("AN21").Value = 0
$ws.Range("AO21").Value = 0
